# API: Gameweeks import (#25)
# Adds two new columns to the "Challenges" import template:
#   - "Show Statistics Continuously" (boolean-as-text flag)
#   - "Gameweek" (numeric)

$wb = $excel.ActiveWorkbook
$wsChallenges = $wb.Worksheets.Item("Challenges")

# New header cells (row 1)
$wsChallenges.Range("S1").Value = "Show Statistics Continuously"
$wsChallenges.Range("T1").Value = "Gameweek"

# New sample-data cells (row 2).
# The "Show Statistics Continuously" example value must be stored as the
# literal text "true" (not an auto-converted boolean), so build it as a
# text formula and flatten it to a plain value in place.
$wsChallenges.Range("S2").Formula = "=""true"""
$wsChallenges.Range("S2").Copy()
$wsChallenges.Range("S2").PasteSpecial(-4163)

$wsChallenges.Range("T2").Value = 1
